$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-18 12:15:39"
$wsZhCn.Range("E5").Value = "2016-03-18 12:15:39"
$wsZhCn.Range("H3").Value = "2016-03-18 12:15:57"
$wsZhCn.Range("H5").Value = "2016-03-18 12:15:57"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-18 12:15:42"
$wsDeDe.Range("E5").Value = "2016-03-18 12:15:42"
$wsDeDe.Range("H3").Value = "2016-03-18 12:16:03"
$wsDeDe.Range("H5").Value = "2016-03-18 12:16:03"
